$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("C3").Value = "-"
$ws.Range("D3").Value = "-"
$ws.Range("E3").Value = "MCT-2A-EAP"
$ws.Range("F3").Value = "-"

# Row 4
$ws.Range("E4").Value = "MCT-2A-EAP"

# Row 6
$ws.Range("C6").Value = "MEC-3A-EAP"

# Row 7
$ws.Range("C7").Value = "MEC-3A-EAP"
